# Insert a new row at position 76 (shifting existing rows 76-151 down to 77-152)
# and populate it with a new weekly price observation for Jengibre.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(76).Insert()

$ws.Cells.Item(76, 1).Value = 9
$ws.Cells.Item(76, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(76, 3).Value = "Metropolitana"
$ws.Cells.Item(76, 4).Value = 45264
$ws.Cells.Item(76, 5).Value = 13
$ws.Cells.Item(76, 6).Value = 100114007
$ws.Cells.Item(76, 7).Value = "Jengibre"
$ws.Cells.Item(76, 8).Value = "Sin especificar"
$ws.Cells.Item(76, 9).Value = "Primera"
$ws.Cells.Item(76, 10).Value = 520
$ws.Cells.Item(76, 11).Value = 15000
$ws.Cells.Item(76, 12).Value = 16000
$ws.Cells.Item(76, 13).Value = 15500
$ws.Cells.Item(76, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(76, 15).Value = "Perú"
$ws.Cells.Item(76, 16).Value = 1192
$ws.Cells.Item(76, 17).Value = 13
$ws.Cells.Item(76, 18).Value = "Hortaliza"
